$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@('253e459e0dc3de445fa4eea651a139ce', 'Mohapatra & Verma', '& & Nona')
    ,@('9d473a7444ef6f4d080d54aec4891d54', 'Balinca & Malaviya', 'Balinca & Malaviya')
    ,@('f35ef701dc13b011bec2afcf484c8012', '& & Nona', 'Jain & Soto')
    ,@('84f2719afb36ce7fdbc3beb60db757e2', 'Balinca & Malaviya', 'Balinca & Malaviya')
    ,@('f7e13e407a5118847f496248141430dc', 'Balinca & Malaviya', '& & Nona')
    ,@('04836fee5898f9fdabbd670dd2653d6f', 'Mohapatra & Verma', 'Mohapatra & Verma')
    ,@('de2cc1e82c2afc223353af87dde50120', 'Mohapatra & Verma', 'Mohapatra & Verma')
    ,@('c2af907400143736d98450762d476166', 'Jonas & Shah', 'Mohapatra & Verma')
    ,@('7b2e92147ae6b6df75feff482e8cd5b6', 'Jain & Soto', 'Jain & Soto')
    ,@('04a58f03b18e2a5bd95f07b690d1e41d', 'Jain & Soto', 'Jain & Soto')
    ,@('04a31d23d25ac871a736c56c72019452', 'Balinca & Malaviya', 'Jonas & Shah')
    ,@('59de1679647fb75c8fafbfe0f4de61c8', 'Jain & Soto', 'Jonas & Shah')
    ,@('b0c7cbbcfa39294c5921ae5efc085d5c', 'Mohapatra & Verma', 'Mohapatra & Verma')
    ,@('2882092cedd64231c4ed11d4cb5e4d97', 'Zhao & jagtap', 'Mohapatra & Verma')
    ,@('91613214ec3950a0599d5a096dc32aa6', 'Lin & Nguyen', 'Diehl & Dowler')
    ,@('e4845b639f3f3d68c9cb38a989fb12da', 'Lin & Nguyen', 'Mohapatra & Verma')
    ,@('b5f9449672616417b46e61686d22006c', 'Diehl & Dowler', 'Kim & Sumesh')
    ,@('546524b567831e12f2a3054e080c6f85', 'Lin & Nguyen', 'Kim & Sumesh')
    ,@('273f10b7459e5be69a257b8acd7a36f0', 'Zhao & jagtap', 'Diehl & Dowler')
    ,@('70dac1be21ceb1ddfba60064a770a883', 'Kim & Sumesh', 'Zhao & jagtap')
    ,@('91b4647cf403ae6f3fc4b0b087947417', 'Mohapatra & Verma', 'Diehl & Dowler')
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $r++
}
